$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 2).Value = 91810  # B5: 91809 -> 91810
$ws.Cells.Item(6, 2).Value = 91810  # B6: 91809 -> 91810
$ws.Cells.Item(7, 2).Value = 79245  # B7: 79244 -> 79245
$ws.Cells.Item(8, 2).Value = 91810  # B8: 91809 -> 91810
$ws.Cells.Item(9, 2).Value = 92181  # B9: 92180 -> 92181
$ws.Cells.Item(10, 2).Value = 92108  # B10: 92107 -> 92108
$ws.Cells.Item(11, 2).Value = 92108  # B11: 92107 -> 92108
$ws.Cells.Item(12, 2).Value = 92023  # B12: 92022 -> 92023
$ws.Cells.Item(13, 1).Value = 131106312  # A13: 131106314 -> 131106312
$ws.Cells.Item(13, 2).Value = 92108  # B13: 91809 -> 92108
$ws.Cells.Item(13, 5).Value = 658  # E13: 1202 -> 658
$ws.Cells.Item(13, 6).Value = "Rosenticka"  # F13: 'Ullticka' -> 'Rosenticka'
$ws.Cells.Item(13, 7).Value = "Fomitopsis rosea"  # G13: 'Phellinidium ferrugineofuscum' -> 'Fomitopsis rosea'
$ws.Cells.Item(13, 8).Value = "(Alb. & Schwein.:Fr.) P.Karst."  # H13: '(P.Karst.) Fiasson & Niemelä' -> '(Alb. & Schwein.:Fr.) P.Karst.'
$ws.Cells.Item(13, 9).Value = ""  # I13: clear (was '1')
$ws.Cells.Item(13, 17).Value = 601540  # Q13: 601556 -> 601540
$ws.Cells.Item(13, 18).Value = 6992576  # R13: 6992605 -> 6992576
$ws.Cells.Item(13, 24).Value = "2025_0872"  # X13: '2025_0870' -> '2025_0872'
$ws.Cells.Item(13, 26).Value = "13:29"  # Z13: '13:21' -> '13:29'
$ws.Cells.Item(13, 28).Value = "13:29"  # AB13: '13:21' -> '13:29'
$ws.Cells.Item(13, 50).Value = "Alexander Hoffmann"  # AX13: 'David Isaksson' -> 'Alexander Hoffmann'
$ws.Cells.Item(14, 1).Value = 131106314  # A14: 131106312 -> 131106314
$ws.Cells.Item(14, 2).Value = 91810  # B14: 92107 -> 91810
$ws.Cells.Item(14, 5).Value = 1202  # E14: 658 -> 1202
$ws.Cells.Item(14, 6).Value = "Ullticka"  # F14: 'Rosenticka' -> 'Ullticka'
$ws.Cells.Item(14, 7).Value = "Phellinidium ferrugineofuscum"  # G14: 'Fomitopsis rosea' -> 'Phellinidium ferrugineofuscum'
$ws.Cells.Item(14, 8).Value = "(P.Karst.) Fiasson & Niemelä"  # H14: '(Alb. & Schwein.:Fr.) P.Karst.' -> '(P.Karst.) Fiasson & Niemelä'
$ws.Cells.Item(14, 9).NumberFormat = "@"
$ws.Cells.Item(14, 9).Value = "1"  # I14: '' -> '1'
$ws.Cells.Item(14, 17).Value = 601556  # Q14: 601540 -> 601556
$ws.Cells.Item(14, 18).Value = 6992605  # R14: 6992576 -> 6992605
$ws.Cells.Item(14, 24).Value = "2025_0870"  # X14: '2025_0872' -> '2025_0870'
$ws.Cells.Item(14, 26).Value = "13:21"  # Z14: '13:29' -> '13:21'
$ws.Cells.Item(14, 28).Value = "13:21"  # AB14: '13:29' -> '13:21'
$ws.Cells.Item(14, 50).Value = "David Isaksson"  # AX14: 'Alexander Hoffmann' -> 'David Isaksson'
$ws.Cells.Item(15, 2).Value = 91810  # B15: 91809 -> 91810
$ws.Cells.Item(16, 2).Value = 91810  # B16: 91809 -> 91810
$ws.Cells.Item(19, 2).Value = 92108  # B19: 92107 -> 92108
$ws.Cells.Item(20, 2).Value = 91810  # B20: 91809 -> 91810
$ws.Cells.Item(21, 2).Value = 91810  # B21: 91809 -> 91810
$ws.Cells.Item(22, 2).Value = 91810  # B22: 91809 -> 91810
$ws.Cells.Item(23, 1).Value = 131106323  # A23: 131106311 -> 131106323
$ws.Cells.Item(23, 2).Value = 92108  # B23: 91809 -> 92108
$ws.Cells.Item(23, 5).Value = 658  # E23: 1202 -> 658
$ws.Cells.Item(23, 6).Value = "Rosenticka"  # F23: 'Ullticka' -> 'Rosenticka'
$ws.Cells.Item(23, 7).Value = "Fomitopsis rosea"  # G23: 'Phellinidium ferrugineofuscum' -> 'Fomitopsis rosea'
$ws.Cells.Item(23, 8).Value = "(Alb. & Schwein.:Fr.) P.Karst."  # H23: '(P.Karst.) Fiasson & Niemelä' -> '(Alb. & Schwein.:Fr.) P.Karst.'
$ws.Cells.Item(23, 17).Value = 601607  # Q23: 601498 -> 601607
$ws.Cells.Item(23, 18).Value = 6992738  # R23: 6992583 -> 6992738
$ws.Cells.Item(23, 24).Value = "2025_0860"  # X23: '2025_0873' -> '2025_0860'
$ws.Cells.Item(23, 26).Value = "12:35"  # Z23: '13:32' -> '12:35'
$ws.Cells.Item(23, 28).Value = "12:35"  # AB23: '13:32' -> '12:35'
$ws.Cells.Item(23, 50).Value = "David Isaksson"  # AX23: 'Alexander Hoffmann' -> 'David Isaksson'
$ws.Cells.Item(24, 1).Value = 131106311  # A24: 131106327 -> 131106311
$ws.Cells.Item(24, 2).Value = 91810  # B24: 91809 -> 91810
$ws.Cells.Item(24, 10).Value = ""  # J24: clear (was 'mycel')
$ws.Cells.Item(24, 17).Value = 601498  # Q24: 601607 -> 601498
$ws.Cells.Item(24, 18).Value = 6992583  # R24: 6992789 -> 6992583
$ws.Cells.Item(24, 24).Value = "2025_0873"  # X24: '2025_0856' -> '2025_0873'
$ws.Cells.Item(24, 26).Value = "13:32"  # Z24: '12:10' -> '13:32'
$ws.Cells.Item(24, 28).Value = "13:32"  # AB24: '12:10' -> '13:32'
$ws.Cells.Item(24, 50).Value = "Alexander Hoffmann"  # AX24: 'David Isaksson, Alexander Hoffmann' -> 'Alexander Hoffmann'
$ws.Cells.Item(25, 1).Value = 131106327  # A25: 131106323 -> 131106327
$ws.Cells.Item(25, 2).Value = 91810  # B25: 92107 -> 91810
$ws.Cells.Item(25, 5).Value = 1202  # E25: 658 -> 1202
$ws.Cells.Item(25, 6).Value = "Ullticka"  # F25: 'Rosenticka' -> 'Ullticka'
$ws.Cells.Item(25, 7).Value = "Phellinidium ferrugineofuscum"  # G25: 'Fomitopsis rosea' -> 'Phellinidium ferrugineofuscum'
$ws.Cells.Item(25, 8).Value = "(P.Karst.) Fiasson & Niemelä"  # H25: '(Alb. & Schwein.:Fr.) P.Karst.' -> '(P.Karst.) Fiasson & Niemelä'
$ws.Cells.Item(25, 10).Value = "mycel"  # J25: None -> 'mycel'
$ws.Cells.Item(25, 18).Value = 6992789  # R25: 6992738 -> 6992789
$ws.Cells.Item(25, 24).Value = "2025_0856"  # X25: '2025_0860' -> '2025_0856'
$ws.Cells.Item(25, 26).Value = "12:10"  # Z25: '12:35' -> '12:10'
$ws.Cells.Item(25, 28).Value = "12:10"  # AB25: '12:35' -> '12:10'
$ws.Cells.Item(25, 50).Value = "David Isaksson, Alexander Hoffmann"  # AX25: 'David Isaksson' -> 'David Isaksson, Alexander Hoffmann'
$ws.Cells.Item(27, 2).Value = 92108  # B27: 92107 -> 92108
$ws.Cells.Item(28, 2).Value = 92269  # B28: 92268 -> 92269
